$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.924.03"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.004.26"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("E11").Value = "  -5.33%  "
$ws.Range("D12").Value = "2.302.51"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("D17").Value = "1.986.64"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "36.832.34"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "0.0$([char]0x2083)0810"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "221.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  -6.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("E30").Value = "  -6.23%  "
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  -5.40%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -5.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "1.451.20"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0210"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0907"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.98%  "
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.993"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "2.192.99"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  -10.47%  "
